$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D19: drop the "treatment satisfaction [OGMS:0000090]; " segment, keeping
# the leading zero-width no-break space (U+FEFF) directly after the ';'
$ws.Range("D19").Value = "injury [OGMS:0000102];" + [char]0xFEFF + "disease [OGMS:0000031]"

# New row 28: GSSO import entry
$ws.Range("A28").Value = "GSSO"
$ws.Range("B28").Value = "http://purl.obolibrary.org/obo/gsso.owl"
$ws.Range("C28").Value = "entity [BFO:0000001]"
$ws.Range("D28").Value = "advocacy organisation [GSSO:005379]; health organisation [GSSO:007328]; human rights organisation [GSSO:003501]; non-profit organisation [GSSO:004615]; money [GSSO:010609]"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "all"
